$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.49809999999998
$ws.Range("A6").Value = -22.70900000000001
$ws.Range("A7").Value = -19.27749999999999
$ws.Range("C7").Value = -12.12890000000001
$ws.Range("C12").Value = -10.5584
$ws.Range("E13").Value = 16.69790000000001
$ws.Range("E14").Value = 16.91990000000001
$ws.Range("C15").Value = -14.9144
$ws.Range("A16").Value = -21.93520000000001
$ws.Range("E16").Value = 16.30790000000001
$ws.Range("E19").Value = 16.45269999999999
$ws.Range("A20").Value = -19.09089999999999
$ws.Range("C20").Value = -11.6449
$ws.Range("C21").Value = -11.7429
$ws.Range("C22").Value = -11.9326
$ws.Range("E22").Value = 17.12650000000002
$ws.Range("C23").Value = -12.29920000000001
$ws.Range("A28").Value = -21.71819999999999
$ws.Range("A29").Value = -21.67349999999998
$ws.Range("C29").Value = -11.8178
$ws.Range("A32").Value = -21.14719999999999
$ws.Range("C34").Value = -11.16290000000001
$ws.Range("E36").Value = 16.2862
$ws.Range("A40").Value = -20.37130000000001
$ws.Range("C42").Value = -12.0422
$ws.Range("C43").Value = -12.90499999999999
$ws.Range("C44").Value = -14.13379999999999
$ws.Range("C45").Value = -14.07719999999999
$ws.Range("A46").Value = -22.0038
$ws.Range("C46").Value = -12.86349999999999
$ws.Range("E46").Value = 16.9687
$ws.Range("C50").Value = -14.23119999999999
$ws.Range("E50").Value = 16.40029999999999
$ws.Range("A51").Value = -22.0234
$ws.Range("C51").Value = -12.2545
$ws.Range("A52").Value = -22.18
$ws.Range("A57").Value = -22.62450000000003
$ws.Range("A59").Value = -22.10120000000001
$ws.Range("A62").Value = -22.04190000000001
$ws.Range("A66").Value = -22.0036
$ws.Range("C66").Value = -12.3541
$ws.Range("C67").Value = -10.8029
$ws.Range("A73").Value = -20.53040000000001
$ws.Range("A74").Value = -21.96059999999998
$ws.Range("C79").Value = -11.925
$ws.Range("C84").Value = -14.07999999999999
$ws.Range("A92").Value = -21.75409999999998
$ws.Range("C92").Value = -11.01270000000001
$ws.Range("E95").Value = 18.13640000000002
$ws.Range("C97").Value = -11.3484
$ws.Range("E97").Value = 16.51959999999999
$ws.Range("A100").Value = -22.07989999999999
